$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("B2").Value = 14.26366496892708
$ws.Range("C2").Value = 11.10383122448297
$ws.Range("D2").Value = 6.826510824745437
$ws.Range("E2").Value = 12.73965231580043
$ws.Range("F2").Value = 44.35775937862943
$ws.Range("I2").Value = 29.52759519459213
$ws.Range("J2").Value = 10.37009018898852
$ws.Range("K2").Value = 15.3236777357165
$ws.Range("N2").Value = 21.97709096228463
$ws.Range("B3").Value = 14.06300550304535
$ws.Range("C3").Value = 10.95478535355906
$ws.Range("D3").Value = 6.799054472679343
$ws.Range("E3").Value = 12.69667285973817
$ws.Range("F3").Value = 44.30839309715767
$ws.Range("I3").Value = 29.56263075874984
$ws.Range("J3").Value = 10.37806476201736
$ws.Range("K3").Value = 15.19114164693388
$ws.Range("N3").Value = 22.03402453523888
$ws.Range("B4").Value = 13.94228964832556
$ws.Range("C4").Value = 10.86551881285995
$ws.Range("D4").Value = 6.783485936586561
$ws.Range("E4").Value = 12.67304878075018
$ws.Range("F4").Value = 44.28792935998142
$ws.Range("I4").Value = 29.5895240158067
$ws.Range("J4").Value = 10.38465996792817
$ws.Range("K4").Value = 15.1130951457312
$ws.Range("N4").Value = 22.07090824167022
$ws.Range("B5").Value = 13.89378679102743
$ws.Range("C5").Value = 10.82975257494245
$ws.Range("D5").Value = 6.777470976452641
$ws.Range("E5").Value = 12.66412392969246
$ws.Range("F5").Value = 44.28206962835113
$ws.Range("I5").Value = 29.6018338399651
$ws.Range("J5").Value = 10.38777472127825
$ws.Range("K5").Value = 15.08216023924443
$ws.Range("N5").Value = 22.08642363047205
$ws.Range("B6").Value = 13.88577643082105
$ws.Range("C6").Value = 10.82385179452109
$ws.Range("D6").Value = 6.776492227696274
$ws.Range("E6").Value = 12.66268455484085
$ws.Range("F6").Value = 44.28124643240006
$ws.Range("I6").Value = 29.60395936863678
$ws.Range("J6").Value = 10.38831771997561
$ws.Range("K6").Value = 15.07707694989734
$ws.Range("N6").Value = 22.08902925854161
$ws.Range("B7").Value = 13.94163264545358
$ws.Range("C7").Value = 10.86503392572701
$ws.Range("D7").Value = 6.783403476871785
$ws.Range("E7").Value = 12.67292556576685
$ws.Range("F7").Value = 44.28784029180577
$ws.Range("I7").Value = 29.58968456557144
$ws.Range("J7").Value = 10.38470024516433
$ws.Range("K7").Value = 15.11267438447519
$ws.Range("N7").Value = 22.07111552320031
$ws.Range("B8").Value = 14.19399939274691
$ws.Range("C8").Value = 11.05200111653325
$ws.Range("D8").Value = 6.816779307127001
$ws.Range("E8").Value = 12.72426331083819
$ws.Range("F8").Value = 44.33869651950496
$ws.Range("I8").Value = 29.53855712255336
$ws.Range("J8").Value = 10.3724872563529
$ws.Range("K8").Value = 15.27730762393885
$ws.Range("N8").Value = 21.99632222618154
$ws.Range("B9").Value = 14.70580014342941
$ws.Range("C9").Value = 11.43444057971949
$ws.Range("D9").Value = 6.892225652731784
$ws.Range("E9").Value = 12.84654757021664
$ws.Range("F9").Value = 44.51634485075315
$ws.Range("I9").Value = 29.48109831075304
$ws.Range("J9").Value = 10.36201449057142
$ws.Range("K9").Value = 15.62510029526663
$ws.Range("N9").Value = 21.8649087720981
$ws.Range("B10").Value = 15.08829679470493
$ws.Range("C10").Value = 11.72227335233791
$ws.Range("D10").Value = 6.953409199125382
$ws.Range("E10").Value = 12.9490708460462
$ws.Range("F10").Value = 44.69396872344506
$ws.Range("I10").Value = 29.46510269913129
$ws.Range("J10").Value = 10.36252830157745
$ws.Range("K10").Value = 15.89375124522731
$ws.Range("N10").Value = 21.77761988052255
$ws.Range("B11").Value = 15.26290588042945
$ws.Range("C11").Value = 11.85411381437425
$ws.Range("D11").Value = 6.982413004547946
$ws.Range("E11").Value = 12.99834462361378
$ws.Range("F11").Value = 44.78487997421046
$ws.Range("I11").Value = 29.46353892312532
$ws.Range("J11").Value = 10.36454038056826
$ws.Range("K11").Value = 16.01835470078392
$ws.Range("N11").Value = 21.73991257929062
$ws.Range("B12").Value = 15.3290478052836
$ws.Range("C12").Value = 11.9041197122912
$ws.Range("D12").Value = 6.993557243679066
$ws.Range("E12").Value = 13.01737107500129
$ws.Range("F12").Value = 44.82074577326816
$ws.Range("I12").Value = 29.46376911143746
$ws.Range("J12").Value = 10.36555744079664
$ws.Range("K12").Value = 16.06584309148691
$ws.Range("N12").Value = 21.72592099939128
$ws.Range("B13").Value = 15.31480305093265
$ws.Range("C13").Value = 11.89334722738046
$ws.Range("D13").Value = 6.991150084711633
$ws.Range("E13").Value = 13.01325723026878
$ws.Range("F13").Value = 44.81295764432853
$ws.Range("I13").Value = 29.46368295213026
$ws.Range("J13").Value = 10.3653270625609
$ws.Range("K13").Value = 16.05560275127105
$ws.Range("N13").Value = 21.72892156171945
$ws.Range("B14").Value = 15.26834736432949
$ws.Range("C14").Value = 11.85822648901698
$ws.Range("D14").Value = 6.983326661140728
$ws.Range("E14").Value = 12.99990264279201
$ws.Range("F14").Value = 44.78780190859443
$ws.Range("I14").Value = 29.46354137605016
$ws.Range("J14").Value = 10.36461894541174
$ws.Range("K14").Value = 16.02225574314954
$ws.Range("N14").Value = 21.73875572653503
$ws.Range("B15").Value = 15.23989274279122
$ws.Range("C15").Value = 11.83672309905002
$ws.Range("D15").Value = 6.9785553516167
$ws.Range("E15").Value = 12.99177008735628
$ws.Range("F15").Value = 44.77258033548246
$ws.Range("I15").Value = 29.4635617706988
$ws.Range("J15").Value = 10.36421840877166
$ws.Range("K15").Value = 16.00186809618307
$ws.Range("N15").Value = 21.7448168483788
$ws.Range("B16").Value = 15.07689222998696
$ws.Range("C16").Value = 11.71367122875358
$ws.Range("D16").Value = 6.951536681746997
$ws.Range("E16").Value = 12.94590275747106
$ws.Range("F16").Value = 44.68822972769214
$ws.Range("I16").Value = 29.46531993403201
$ws.Range("J16").Value = 10.3624325444787
$ws.Range("K16").Value = 15.88565282609643
$ws.Range("N16").Value = 21.78012437199717
$ws.Range("B17").Value = 14.9770050248034
$ws.Range("C17").Value = 11.63837909858239
$ws.Range("D17").Value = 6.935256456710444
$ws.Range("E17").Value = 12.9184316992592
$ws.Range("F17").Value = 44.63906320390273
$ws.Range("I17").Value = 29.46786242424812
$ws.Range("J17").Value = 10.36179206946061
$ws.Range("K17").Value = 15.81494279404495
$ws.Range("N17").Value = 21.80229659451929
$ws.Range("B18").Value = 14.91961378587485
$ws.Range("C18").Value = 11.59516102893654
$ws.Range("D18").Value = 6.926003184418293
$ws.Range("E18").Value = 12.90287985714585
$ws.Range("F18").Value = 44.61173640912274
$ws.Range("I18").Value = 29.4698624774727
$ws.Range("J18").Value = 10.36159111503053
$ws.Range("K18").Value = 15.77449991463625
$ws.Range("N18").Value = 21.81523782035021
$ws.Range("B19").Value = 14.9001946669128
$ws.Range("C19").Value = 11.58054475710715
$ws.Range("D19").Value = 6.922889412683478
$ws.Range("E19").Value = 12.89765733750454
$ws.Range("F19").Value = 44.60264801366792
$ws.Range("I19").Value = 29.47063197494052
$ws.Range("J19").Value = 10.3615518475388
$ws.Range("K19").Value = 15.76084697477222
$ws.Range("N19").Value = 21.81965185784673
$ws.Range("B20").Value = 14.98763233659731
$ws.Range("C20").Value = 11.64638532908092
$ws.Range("D20").Value = 6.936978111127691
$ws.Range("E20").Value = 12.92133037292773
$ws.Range("F20").Value = 44.64419858076928
$ws.Range("I20").Value = 29.46753611623681
$ws.Range("J20").Value = 10.36184292464851
$ws.Range("K20").Value = 15.82244672671514
$ws.Range("N20").Value = 21.79991683273945
$ws.Range("B21").Value = 15.28199246895727
$ws.Range("C21").Value = 11.86854049603316
$ws.Range("D21").Value = 6.985620276562425
$ws.Range("E21").Value = 13.00381532801858
$ws.Range("F21").Value = 44.79515180629823
$ws.Range("I21").Value = 29.46356063690298
$ws.Range("J21").Value = 10.36482001781666
$ws.Range("K21").Value = 16.03204264617884
$ws.Range("N21").Value = 21.73585939794237
$ws.Range("B22").Value = 15.47446825684585
$ws.Range("C22").Value = 12.01418052415566
$ws.Range("D22").Value = 7.018346372221512
$ws.Range("E22").Value = 13.05986060285647
$ws.Range("F22").Value = 44.902192233322
$ws.Range("I22").Value = 29.46575586999021
$ws.Range("J22").Value = 10.36825244356257
$ws.Range("K22").Value = 16.17077745272054
$ws.Range("N22").Value = 21.69566894449095
$ws.Range("B23").Value = 15.3717531195346
$ws.Range("C23").Value = 11.93642461421895
$ws.Range("D23").Value = 7.000796696044026
$ws.Range("E23").Value = 13.02975661866392
$ws.Range("F23").Value = 44.84430073418741
$ws.Range("I23").Value = 29.46414544946425
$ws.Range("J23").Value = 10.36628468784254
$ws.Range("K23").Value = 16.09658516160248
$ws.Range("N23").Value = 21.71696621708448
$ws.Range("B24").Value = 14.98282761402689
$ws.Range("C24").Value = 11.64276549550993
$ws.Range("D24").Value = 6.936199419063168
$ws.Range("E24").Value = 12.92001912905573
$ws.Range("F24").Value = 44.64187394859035
$ws.Range("I24").Value = 29.46768196312969
$ws.Range("J24").Value = 10.36181941196347
$ws.Range("K24").Value = 15.81905354382527
$ws.Range("N24").Value = 21.80099211825756
$ws.Range("B25").Value = 14.56593065830297
$ws.Range("C25").Value = 11.32957514812424
$ws.Range("D25").Value = 6.870779082263994
$ws.Range("E25").Value = 12.8112008202844
$ws.Range("F25").Value = 44.45997784860936
$ws.Range("I25").Value = 29.49204571029573
$ws.Range("J25").Value = 10.36340517480075
$ws.Range("K25").Value = 15.52855282792171
$ws.Range("N25").Value = 21.89883054336979
